$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that need to stay as TEXT even though the new value looks numeric:
# mark them with a temporary text number format, then reset style afterward
$textCells = @("D5", "D6", "D8", "D10", "D11", "D12", "D13", "D14", "D17", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D30", "D31", "D32", "D33", "D36", "D37", "D38", "D39", "D40", "D42", "D43", "D44", "D45", "D46", "D47", "D50")
foreach ($addr in $textCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = '60.484.41'
$ws.Range("E2").Value = '  -2.29%  '
$ws.Range("D3").Value = '2.897.44'
$ws.Range("E3").Value = '  -3.39%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '583.13'
$ws.Range("E5").Value = '  -1.61%  '
$ws.Range("D6").Value = '146.62'
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '0.501'
$ws.Range("E8").Value = '  -2.73%  '
$ws.Range("D9").Value = '2.889.76'
$ws.Range("E9").Value = '  -3.57%  '
$ws.Range("D10").Value = '6.68'
$ws.Range("E10").Value = '  +7.57%  '
$ws.Range("D11").Value = '0.143'
$ws.Range("E11").Value = '  -3.37%  '
$ws.Range("D12").Value = '0.444'
$ws.Range("E12").Value = '  -2.61%  '
$ws.Range("D13").Value = '0.0000223'
$ws.Range("E13").Value = '  -3.31%  '
$ws.Range("D14").Value = '34.21'
$ws.Range("E14").Value = '  -0.01%  '
$ws.Range("E15").Value = '  +0.00%  '
$ws.Range("D16").Value = '3.380.28'
$ws.Range("E16").Value = '  -3.29%  '
$ws.Range("B17").Value = 'Polkadot'
$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D17").Value = '6.80'
$ws.Range("E17").Value = '  -2.44%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '60.468.92'
$ws.Range("E18").Value = '  -2.28%  '
$ws.Range("D19").Value = '2.902.68'
$ws.Range("E19").Value = '  -3.27%  '
$ws.Range("D20").Value = '423.46'
$ws.Range("E20").Value = '  -5.23%  '
$ws.Range("D21").Value = '13.55'
$ws.Range("E21").Value = '  -3.99%  '
$ws.Range("D22").Value = '0.667'
$ws.Range("E22").Value = '  -2.43%  '
$ws.Range("D23").Value = '7.11'
$ws.Range("E23").Value = '  -3.35%  '
$ws.Range("D24").Value = '80.68'
$ws.Range("E24").Value = '  -1.82%  '
$ws.Range("D25").Value = '10.90'
$ws.Range("E25").Value = '  -0.10%  '
$ws.Range("D26").Value = '2.16'
$ws.Range("E26").Value = '  -2.53%  '
$ws.Range("D27").Value = '11.77'
$ws.Range("E27").Value = '  -2.77%  '
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("D30").Value = '7.23'
$ws.Range("E30").Value = '  +0.40%  '
$ws.Range("D31").Value = '2.17'
$ws.Range("E31").Value = '  +4.34%  '
$ws.Range("D32").Value = '2.61'
$ws.Range("E32").Value = '  -3.27%  '
$ws.Range("D33").Value = '26.55'
$ws.Range("E33").Value = '  -3.13%  '
$ws.Range("E34").Value = '  -4.36%  '
$ws.Range("D35").Value = '0.0₃0833'
$ws.Range("E35").Value = '  -0.87%  '
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  -1.83%  '
$ws.Range("D37").Value = '5.64'
$ws.Range("E37").Value = '  -2.77%  '
$ws.Range("D38").Value = '49.62'
$ws.Range("E38").Value = '  -1.04%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").Value = '2.02'
$ws.Range("E39").Value = '  -0.36%  '
$ws.Range("B40").Value = 'dogwifhat'
$ws.Range("C40").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D40").Value = '2.94'
$ws.Range("E40").Value = '  +1.29%  '
$ws.Range("E41").Value = '  +0.48%  '
$ws.Range("D42").Value = '8.71'
$ws.Range("E42").Value = '  -2.79%  '
$ws.Range("D43").Value = '0.288'
$ws.Range("E43").Value = '  +3.43%  '
$ws.Range("D44").Value = '41.07'
$ws.Range("E44").Value = '  +1.15%  '
$ws.Range("D45").Value = '371.20'
$ws.Range("E45").Value = '  -6.11%  '
$ws.Range("D46").Value = '0.0344'
$ws.Range("E46").Value = '  -1.66%  '
$ws.Range("B47").Value = 'Monero'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D47").Value = '132.79'
$ws.Range("E47").Value = '  +0.20%  '
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '2.645.03'
$ws.Range("E48").Value = '  -2.51%  '
$ws.Range("E49").Value = '  +0.00%  '
$ws.Range("D50").Value = '25.15'
$ws.Range("E50").Value = '  +6.31%  '
$ws.Range("E51").Value = '  -0.88%  '

# Restore default style on the text-forced numeric cells (keeps them as text,
# drops the temporary NumberFormat so no stray formatting is introduced)
foreach ($addr in $textCells) { $ws.Range($addr).Style = "Normal" }
